# Updates the exception stack trace text in the document to reflect
# moving the test/build environment from Java 8 (3.2.1) to a newer JDK (3.2.2):
#  - "sun.reflect.*" / "java.lang.reflect.Method" frames become "java.base/..." module-qualified frames
#  - several line numbers in referenced source files shift
#  - the JUnit/Eclipse test runner trace is replaced by a Maven Surefire / Tycho / Equinox launcher trace
$d = $word.ActiveDocument

$old1 = "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)"
$new1 = "`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:564)"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) { throw "Replacement 1 did not find its target text" }
Write-Host "Replace 1: OK"

$old2 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)"
$new2 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:186)"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) { throw "Replacement 2 did not find its target text" }
Write-Host "Replace 2: OK"

$old3 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)"
$new3 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)"
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) { throw "Replacement 3 did not find its target text" }
Write-Host "Replace 3: OK"

$old4 = "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)`n`tat sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)"
$new4 = "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:529)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:415)`n`tat jdk.internal.reflect.GeneratedMethodAccessor10.invoke(Unknown Source)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:564)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)"
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $found4) { throw "Replacement 4 did not find its target text" }
Write-Host "Replace 4: OK"

$old5 = "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)"
$new5 = "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)"
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $found5) { throw "Replacement 5 did not find its target text" }
Write-Host "Replace 5: OK"

$old6 = "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$new6 = "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $found6) { throw "Replacement 6 did not find its target text" }
Write-Host "Replace 6: OK"

$old7 = "`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)"
$new7 = "`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)"
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
if (-not $found7) { throw "Replacement 7 did not find its target text" }
Write-Host "Replace 7: OK"

$old8 = "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$new8 = "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
if (-not $found8) { throw "Replacement 8 did not find its target text" }
Write-Host "Replace 8: OK"

$old9 = "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$new9 = "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:564)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:564)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:564)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)"
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
if (-not $found9) { throw "Replacement 9 did not find its target text" }
Write-Host "Replace 9: OK"
